$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.943.98'
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").Value = '3.120.97'
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'531.82"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").Value = "'138.29"
$ws.Range("E6").Value = '  -1.77%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = "'0.463"
$ws.Range("E8").Value = '  +3.31%  '

$ws.Range("E9").Value = '  +0.72%  '

$ws.Range("E10").Value = '  -2.17%  '

$ws.Range("D11").Value = "'0.408"
$ws.Range("E11").Value = '  +1.89%  '

$ws.Range("D12").Value = '3.657.68'
$ws.Range("E12").Value = '  -0.96%  '

$ws.Range("E13").Value = '  +1.41%  '

$ws.Range("D14").Value = "'25.47"
$ws.Range("E14").Value = '  -0.75%  '

$ws.Range("D15").Value = "'0.0000162"
$ws.Range("E15").Value = '  -2.03%  '

$ws.Range("D16").Value = '57.929.42'
$ws.Range("E16").Value = '  -0.16%  '

$ws.Range("D17").Value = '3.120.63'
$ws.Range("E17").Value = '  -0.77%  '

$ws.Range("E18").Value = '  -2.39%  '

$ws.Range("E19").Value = '  -2.23%  '

$ws.Range("D20").Value = "'8.00"
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").Value = "'350.65"
$ws.Range("E21").Value = '  -1.36%  '

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.19%  '

$ws.Range("D23").Value = "'68.95"
$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("D24").Value = "'0.504"
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("E25").Value = '  -1.85%  '

$ws.Range("E26").Value = '  -0.20%  '

$ws.Range("D27").Value = '0.0₃0872'
$ws.Range("E27").Value = '  -7.33%  '

$ws.Range("E28").Value = '  -3.48%  '

$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").Value = "'6.04"
$ws.Range("E30").Value = '  -5.57%  '

$ws.Range("D31").Value = "'21.21"
$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("D32").Value = "'4.92"
$ws.Range("E32").Value = '  +0.88%  '

$ws.Range("E33").Value = '  -5.53%  '

$ws.Range("D34").Value = "'158.64"
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("D35").Value = "'6.04"
$ws.Range("E35").Value = '  -2.70%  '

$ws.Range("D36").Value = "'25.88"
$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("E37").Value = '  -2.67%  '

$ws.Range("D38").Value = "'1.67"
$ws.Range("E38").Value = '  +3.29%  '

$ws.Range("E39").Value = '  -0.48%  '

$ws.Range("E40").Value = '  -2.09%  '

$ws.Range("E41").Value = '  -1.30%  '

$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'37.13"
$ws.Range("E42").Value = '  +1.11%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.393.87'
$ws.Range("E43").Value = '  +2.46%  '

$ws.Range("D44").Value = '3.160.21'
$ws.Range("E44").Value = '  -1.03%  '

$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = '  -0.12%  '

$ws.Range("E46").Value = '  -2.65%  '

$ws.Range("D47").Value = "'0.963"
$ws.Range("E47").Value = '  -4.40%  '

$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("D49").Value = "'19.76"
$ws.Range("E49").Value = '  -3.15%  '

$ws.Range("D50").Value = "'0.739"
$ws.Range("E50").Value = '  -2.77%  '

$ws.Range("E51").Value = '  +1.57%  '
